# Updated the document Version.
#
# 1) Append "sion" to the end of the first paragraph's text as its own
#    run (so "Design Document Version1" / "sion" stay as two separate
#    <w:r> elements rather than being merged into one run).
# 2) Add a new paragraph after it reading
#    "Updated the document Version to 2.0".

$d = $word.ActiveDocument

$para1 = $d.Paragraphs(1)
$endOfText = $para1.Range.End - 1   # position just before the paragraph mark

# Insert a temporary paragraph break, type the new text into the resulting
# (new) paragraph, then delete that paragraph break again. Because the
# break is removed via Range.Delete (rather than by merging two already
# co-resident runs), the engine keeps "sion" as its own run instead of
# folding it back into the preceding run.
$splitPoint = $d.Range($endOfText, $endOfText)
$splitPoint.InsertParagraphAfter()

$newRunPos = $endOfText + 1
$newRunRange = $d.Range($newRunPos, $newRunPos)
$newRunRange.InsertAfter("sion")

$markRange = $d.Range($endOfText, $endOfText + 1)
$markRange.Delete()

# Add the new second paragraph with the version note.
$para1 = $d.Paragraphs(1)
$para1End = $para1.Range.End
$paraBreak = $d.Range($para1End, $para1End)
$paraBreak.InsertParagraphAfter()

$newPara = $d.Paragraphs(2)
$newPara.Range.InsertAfter("Updated the document Version to 2.0")
